# Fruta / hortaliza, semanal
# Insert two new weekly rows (row 20 and 21) above the existing data,
# pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 20 (shifts old rows 20.. down to 22..)
$ws.Rows("20:21").Insert()

# --- New row 20 ---
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Vega Monumental Concepción"
$ws.Range("C20").Value = "Bíobío"
$ws.Range("D20").Value = 44497
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100102
$ws.Range("H20").Value = "Cítricos"
$ws.Range("I20").Value = 100102004
$ws.Range("J20").Value = "Mandarina"
$ws.Range("K20").Value = "Murcott"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 6000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 6500
$ws.Range("Q20").Value = "`$/bandeja 10 kilos"
$ws.Range("R20").Value = "Provincia de Limarí"
$ws.Range("S20").Value = 650
$ws.Range("T20").Value = 10

# --- New row 21 ---
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44497
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100102
$ws.Range("H21").Value = "Cítricos"
$ws.Range("I21").Value = 100102004
$ws.Range("J21").Value = "Mandarina"
$ws.Range("K21").Value = "Murcott"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 5000
$ws.Range("O21").Value = 5000
$ws.Range("P21").Value = 5000
$ws.Range("Q21").Value = "`$/bandeja 10 kilos"
$ws.Range("R21").Value = "Provincia de Limarí"
$ws.Range("S21").Value = 500
$ws.Range("T21").Value = 10

Write-Output "done"
